$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")
$ws.Range("N244:O248").NumberFormat = "General"
Write-Output "done"
